$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Merge Sort Recurrsion"
$ws.Range("C5").Value = "Recurrsion"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Quick Sort Recurrsion"
$ws.Range("C6").Value = "Recurrsion"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Making Subset of the String given by the User "
$ws.Range("C7").Value = "Recurrsion"

$ws.Range("A5:C7").HorizontalAlignment = -4108

$ws.Range("C8").Select()
